$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.965.44"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "1.640.33"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  +0.60%  "

$ws.Range("D5").Value = "216.11"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").Value = "19.64"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E13").Value = "  +1.40%  "

$ws.Range("D14").Value = "1.644.45"
$ws.Range("E14").Value = "  +1.80%  "

$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").Value = "62.87"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "25.926.65"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").Value = "193.05"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  +0.64%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("E24").Value = "  +7.28%  "

$ws.Range("E25").Value = "  +1.33%  "

$ws.Range("D26").Value = "144.51"
$ws.Range("E26").Value = "  +2.22%  "

$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").Value = "6.93"
$ws.Range("E28").Value = "  +1.51%  "

$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("E32").Value = "  -1.01%  "

$ws.Range("E33").Value = "  +1.68%  "

$ws.Range("E34").Value = "  -2.43%  "

$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").Value = "1.133.66"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "0.541"
$ws.Range("E38").Value = "  -0.76%  "

$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("E41").Value = "  +1.87%  "

$ws.Range("D42").Value = "99.39"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "0.794"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").Value = "1.779.75"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("D46").Value = "56.67"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  +0.71%  "
